# Update the "QuantitativeMetrics" sheet of UC3.4_TC1.xlsx to reflect
# the refreshed evaluation run:
#   - C5: shorten the "compilation failed" note
#   - B12 / C12: refreshed CodeBLEU score (and its breakdown note)

$ws = $excel.Worksheets.Item("QuantitativeMetrics")

# C5 - shorter note about the non-existing method call
$ws.Range("C5").Value = "Calling a not existing method"

# B12 - new CodeBLEU score
$ws.Range("B12").Value = 0.2491919426004884

# C12 - breakdown note updated to match the new CodeBLEU score
$ws.Range("C12").Value = "{'codebleu': 0.2491919426004884, 'ngram_match_score': 0.09578916373048774, 'weighted_ngram_match_score': 0.11393541730269516, 'syntax_match_score': 0.5428571428571428, 'dataflow_match_score': 0.2441860465116279}"
